$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: existing (mostly blank) row gets a full entry ---
$ws.Cells.Item(11, 1).Value = 46039
$ws.Cells.Item(11, 2).Value = 0.791666666666667
$ws.Cells.Item(11, 3).Value = 0.854166666666667
$ws.Cells.Item(11, 4).Formula = "=C11-B11"
$ws.Cells.Item(11, 5).Value = "Front"
$ws.Cells.Item(11, 6).Value = "Laatta grafiikat"

# --- Row 12: brand-new row ---
# Copy formatting from the row above first so the new row matches the
# table's existing look (date / time / duration formats).
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)

$ws.Cells.Item(12, 1).Value = 46040
$ws.Cells.Item(12, 2).Value = 0.895833333333333
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 4).Formula = "=C12-B12"
$ws.Cells.Item(12, 5).Value = "Front"
$ws.Cells.Item(12, 6).Value = "Laatta grafiikat"

# --- Row 13: brand-new row ---
$ws.Range("A11:F11").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)

$ws.Cells.Item(13, 1).Value = 46041
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0.0416666666666667
$ws.Cells.Item(13, 4).Formula = "=C13-B13"
$ws.Cells.Item(13, 5).Value = "Front"
$ws.Cells.Item(13, 6).Value = "Laatta grafiikat"

# --- Move the active selection cell (matches the diff's sheetView) ---
$ws.Range("C14").Select()

Write-Output "edit complete"
